# fix(gui) step 1 and 2
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Hoja1")

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the price figures for the two line items (S-005 / T-045)
$ws.Range("D33").Value = 49.59
$ws.Range("D34").Value = 38.976
